$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Projet de site de commerce électronique</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Créez une plateforme de commerce électronique avec Angular pour la partie frontale, ASP.Net Web API pour la partie backend, MySQL pour la base de données et NHibernate ou Fluent Nhibernate pour la couche d'accès aux données. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Ce projet peut inclure des fonctionnalités telles </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">que l'inscription de </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>l'utilisateur</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">, la recherche de </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>produits</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">, l'ajout de produits au </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>panier</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">, la </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>commande</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> de produits, la gestion de la commande, les options de </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>paiement</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> et les </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:u w:val="single"/></w:rPr><w:t>commentaires</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> des clients. </w:t></w:r><w:r><w:t>Ce projet pourrait avoir un niveau de difficulté moyen.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="6096"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Gestion des </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>utilisateurs</w:t></w:r><w:r><w:t xml:space="preserve"> : en plus de </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>l'inscription</w:t></w:r><w:r><w:t xml:space="preserve"> des utilisateurs, vous pouvez également inclure des fonctionnalités de connexion et de déconnexion, de profil utilisateur (avec la possibilité d'ajouter des informations telles que </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>l'adresse</w:t></w:r><w:r><w:t xml:space="preserve"> de livraison et les options de paiement) et de réinitialisation de mot de passe.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestion des </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>produits</w:t></w:r><w:r><w:t xml:space="preserve"> : pour la recherche de produits, vous pouvez ajouter des filtres pour aider les utilisateurs à trouver des produits spécifiques (par exemple, par catégorie, par prix, par marque, par couleur, etc.). Pour l'ajout de produits au panier, vous pouvez inclure une fonctionnalité de quantité et une vérification de disponibilité du produit. Pour la commande de produits, vous pouvez inclure la sélection de la méthode de livraison et le suivi de la commande.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestion des </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>commentaires</w:t></w:r><w:r><w:t xml:space="preserve"> des clients : les utilisateurs peuvent laisser des commentaires sur les produits qu'ils ont achetés. Vous pouvez inclure des fonctionnalités telles que la notation des produits, la modération des commentaires et la réponse aux commentaires.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestion des </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>options de paiement</w:t></w:r><w:r><w:t xml:space="preserve"> : vous pouvez inclure plusieurs options de paiement, telles que les cartes de crédit, les portefeuilles électroniques, les virements bancaires, etc. Vous pouvez également inclure des fonctionnalités de sécurité pour protéger les données de paiement des utilisateurs.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestion des </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>promotions</w:t></w:r><w:r><w:t xml:space="preserve"> : vous pouvez inclure des fonctionnalités de promotion, telles que des codes de réduction, des offres spéciales pour les clients fidèles, etc.</w:t></w:r></w:p><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($newBodyXml)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
